# Record a new game (2025-09-14) into Game_Record, then leave the
# selection/active-sheet state the way the author left it: Game_Record's
# selection sitting one cell past the new row (G12), and Stat_Sheet as the
# active/visible tab with B2 selected.

$wb  = $excel.ActiveWorkbook
$gameRecord = $wb.Worksheets.Item("Game_Record")
$statSheet  = $wb.Worksheets.Item("Stat_Sheet")

# New row of results for game #11.
$gameRecord.Range("A12").Formula = "=ROW()-1"
$gameRecord.Range("B12").Value = 45914
$gameRecord.Range("C12").Value = "Doanage"
$gameRecord.Range("D12").Value = "Player1"
$gameRecord.Range("E12").Value = "DrSystomatix"
$gameRecord.Range("F12").Value = "SimpleJack"

# Leave Game_Record's own selection on G12 (just right of the new data).
$gameRecord.Activate()
$gameRecord.Range("G12").Select()

# Finish on Stat_Sheet (now the active/visible tab) with B2 selected.
$statSheet.Activate()
$statSheet.Range("B2").Select()
